$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/audit-level"
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")
# R5 mirrors the StructureDefinition URL (Fixed Value of Extension.url)
$elem.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/audit-level"
# Z6 is the Binding Value Set URL
$elem.Range("Z6").Value = "https://2rdoc.pt/fhir/ValueSet/audit-levels"
